# Updates cryptos list values (Price / Volume(1h)) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "58.309.35"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.291.69"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "536.79"
$ws.Range("E5").Value = "  -1.82%  "
Set-TextValue "D6" "131.53"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue "D8" "0.589"
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("D9").Value = "2.286.58"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("E12").Value = "  +0.93%  "
Set-TextValue "D13" "0.334"
$ws.Range("E13").Value = "  -0.30%  "
Set-TextValue "D14" "23.70"
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "2.700.23"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "58.251.85"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "2.262.16"
$ws.Range("E18").Value = "  -1.08%  "
Set-TextValue "D19" "10.57"
$ws.Range("E19").Value = "  -0.60%  "
Set-TextValue "D20" "4.20"
Set-TextValue "D21" "314.56"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("E23").Value = "  +0.00%  "
Set-TextValue "D24" "63.30"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("E26").Value = "  -0.07%  "
Set-TextValue "D27" "8.00"
$ws.Range("E27").Value = "  -1.61%  "
Set-TextValue "D28" "1.30"
$ws.Range("E28").Value = "  -0.53%  "
Set-TextValue "D29" "171.08"
$ws.Range("E29").Value = "  +0.28%  "
Set-TextValue "D30" "1.71"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("E31").Value = "  +0.30%  "
Set-TextValue "D32" "1.09"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("E35").Value = "  -0.02%  "
Set-TextValue "D36" "17.86"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -0.83%  "
Set-TextValue "D39" "3.93"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("E40").Value = "  -0.93%  "
Set-TextValue "D41" "290.26"
$ws.Range("E41").Value = "  -3.31%  "
Set-TextValue "D42" "140.27"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("E46").Value = "  -0.24%  "
Set-TextValue "D47" "18.27"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("E48").Value = "  -1.54%  "
Set-TextValue "D49" "10.95"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("E51").Value = "  +1.20%  "
